$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New file identities used throughout the report (this run generated a
# handoff for a newly-renamed source file and updated xliff packages).
# ---------------------------------------------------------------------------
$oldGuid1 = "120e8896-5a44-429b-9c1c-56c44d399a66"
$oldGuid2 = "f198c155-ecc6-4c27-965f-aaeb8c0d8e40"
$newGuid1 = "50665634-df87-4bdf-8844-366d388a1333"
$newGuid2 = "ffff8ec9e819-0f1f-4c14-b7f9-caa692304b6f"

$file1 = "$newGuid1.md"
$file2 = "$newGuid2.md"
$path1 = "e2e\$newGuid1.md"
$path2 = "e2e\$newGuid2.md"

$xlf1zh = "$newGuid1.80542aa4a325e0fa1447904cebd7903e86ddb965.zh-cn.xlf"
$xlf1de = "$newGuid1.80542aa4a325e0fa1447904cebd7903e86ddb965.de-de.xlf"

$status = "Ready for handoff"
$handoffTime = "2016-08-31 05:04:08"
$handoffTimeDeDe = "2016-08-31 05:04:13"
$zeroTime = "0001-01-01 00:00:00"

# URLs backing the hyperlinks - unchanged targets, only their display text
# (the file name) needs to reflect the new guids.
$url1Main = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/57dd73b5422cff2916fcf2efb313ca2be5507b48/e2e/$oldGuid1.md"
$url2Main = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/57dd73b5422cff2916fcf2efb313ca2be5507b48/e2e/$oldGuid2.md"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Hyperlinks.Delete()

$ws1.Range("A2").Value = $file1
$ws1.Range("B2").Value = $path1
$ws1.Range("E2").Value = $status
$ws1.Range("F2").Value = $status
$ws1.Range("G2").Value = $handoffTimeDeDe

$ws1.Range("A3").Value = $file2
$ws1.Range("B3").Value = $path2
$ws1.Range("E3").Value = $status
$ws1.Range("F3").Value = $status
$ws1.Range("G3").Value = $handoffTimeDeDe

$ws1.Hyperlinks.Add($ws1.Range("B2"), $url1Main, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $path1) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B3"), $url2Main, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $path2) | Out-Null

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

$ws2.Range("A2").Value = $file1
$ws2.Range("C2").Value = $status
$ws2.Range("F2").Value = "'True"
$ws2.Range("F2").Style = "Normal"
$ws2.Range("G2").Value = $xlf1zh
$ws2.Range("H2").Value = $handoffTime
$ws2.Range("I2").Value = ""
$ws2.Range("I2").Style = "Normal"
$ws2.Range("J2").Value = ""
$ws2.Range("K2").Value = $zeroTime

$ws2.Range("A3").Value = $file2
$ws2.Range("C3").Value = $status
$ws2.Range("F3").Value = "'True"
$ws2.Range("F3").Style = "Normal"
$ws2.Range("G3").Value = $xlf1de
$ws2.Range("H3").Value = $handoffTime
$ws2.Range("I3").Value = ""
$ws2.Range("I3").Style = "Normal"
$ws2.Range("J3").Value = ""
$ws2.Range("K3").Value = $zeroTime

$ws2.Hyperlinks.Add($ws2.Range("A2"), $url1Main, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $file1) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), $url2Main, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $file2) | Out-Null

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

$ws3.Range("A2").Value = $file1
$ws3.Range("C2").Value = $status
$ws3.Range("G2").Value = $xlf1de
$ws3.Range("H2").Value = $handoffTimeDeDe
$ws3.Range("I2").Value = ""
$ws3.Range("I2").Style = "Normal"
$ws3.Range("J2").Value = ""
$ws3.Range("K2").Value = $zeroTime

$ws3.Range("A3").Value = $file2
$ws3.Range("C3").Value = $status
$ws3.Range("F3").Value = "'True"
$ws3.Range("F3").Style = "Normal"
$ws3.Range("G3").Value = $xlf1de
$ws3.Range("H3").Value = $handoffTimeDeDe
$ws3.Range("I3").Value = ""
$ws3.Range("I3").Style = "Normal"
$ws3.Range("J3").Value = ""
$ws3.Range("K3").Value = $zeroTime

$ws3.Hyperlinks.Add($ws3.Range("A2"), $url1Main, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $file1) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), $url2Main, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $file2) | Out-Null
